$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.467.92"
$ws.Range("D3").Value = "1.913.15"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.73"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4818"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08145"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.45"
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("D12").Value = "1.920.52"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.000"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.128"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.14"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06771"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.68"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "29.487.31"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.185"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").Value = "2.120.78"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.82"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.364"
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.108"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.72"
$ws.Range("E30").Value = "  +1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.024"
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09526"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.515"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.561"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02265"
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06096"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.176"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5965"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.967"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.67"
$ws.Range("E41").Value = "  +5.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1853"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.285"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.399"
$ws.Range("E44").Value = "  -5.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.53"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07632"
$ws.Range("E46").Value = "  -3.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5571"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.71"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.417"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.45"
$ws.Range("E51").Value = "  +1.65%  "
